# Added time series filtration: recompute per-swimmer/per-event scoring
# values on the "Scoring" sheet after filtering out older results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 13
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 12
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("R2").Value = 3
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 7
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 2.25
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 2.25
$ws.Range("R3").Value = 0.75
$ws.Range("S3").Value = 2.5
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 9
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 2.5
$ws.Range("R4").Value = 0.75
$ws.Range("S4").Value = 2.25
$ws.Range("J5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 7
$ws.Range("R7").Value = 0.25
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 0.5
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0.5
$ws.Range("Q8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 9
$ws.Range("N9").Value = 2.75
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.25
$ws.Range("F10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = 13
$ws.Range("J11").Value = 13
$ws.Range("N11").Value = 2.75
$ws.Range("O11").Value = 3.25
$ws.Range("P11").Value = 2.75
$ws.Range("S11").Value = 3.25
$ws.Range("D12").Value = 1
$ws.Range("M12").Value = 0
$ws.Range("B13").Value = 3
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1
$ws.Range("N13").Value = 0.75
$ws.Range("P13").Value = 0.75
